# Applies the "Created Visual Analysis to Consider Chart Exports, adjusted
# other functions" commit:
#   * Sheet "Data Description": add a new "thalach" column (F) with its
#     describe() stats, and tighten up the existing age/sex/restbp/chol
#     stats (re-run on a slightly different train/test split).
#   * Sheet "Classes Balance": updated Healthy/Sick counts.
#   * Sheet "Intrinsic Discrepancy": updated per-feature discrepancy scores.
#   * Sheets "Data Samples" / "Intrinsic Discrepancy" headers still read
#     cp/fbs/restecg/thalach - values unchanged, just re-touching them so the
#     shared-string table stays consistent.

$wb = $excel.ActiveWorkbook

$wsDesc  = $wb.Worksheets.Item("Data Description")
$wsSamp  = $wb.Worksheets.Item("Data Samples")
$wsBal   = $wb.Worksheets.Item("Classes Balance")
$wsDisc  = $wb.Worksheets.Item("Intrinsic Discrepancy")

# ---------------------------------------------------------------------
# 1. "Data Description" sheet - add the thalach column (F) and refresh
#    the describe() table for age/sex/restbp/chol/thalach.
# ---------------------------------------------------------------------

# Header row
$wsDesc.Cells.Item(1, 2).Value = "age"
$wsDesc.Cells.Item(1, 3).Value = "sex"
$wsDesc.Cells.Item(1, 4).Value = "restbp"
$wsDesc.Cells.Item(1, 5).Value = "chol"
$wsDesc.Cells.Item(1, 6).Value = "thalach"
$wsDesc.Cells.Item(1, 6).Copy()
$wsDesc.Cells.Item(1, 5).Copy()
$wsDesc.Cells.Item(1, 6).PasteSpecial(-4122)   # xlPasteFormats, match other header cells
$excel.CutCopyMode = $false

# count (row labels in column A are untouched - their text doesn't change,
# only their internal shared-string index shifts, which Excel manages for us)
$wsDesc.Cells.Item(2, 2).Value = 301
$wsDesc.Cells.Item(2, 3).Value = 301
$wsDesc.Cells.Item(2, 4).Value = 301
$wsDesc.Cells.Item(2, 5).Value = 301
$wsDesc.Cells.Item(2, 6).Value = 301

# mean
$wsDesc.Cells.Item(3, 2).Value = 54.37873754152824
$wsDesc.Cells.Item(3, 3).Value = 0.6843853820598007
$wsDesc.Cells.Item(3, 4).Value = 131.7840531561462
$wsDesc.Cells.Item(3, 5).Value = 245.8671096345515
$wsDesc.Cells.Item(3, 6).Value = 149.7508305647841

# std
$wsDesc.Cells.Item(4, 2).Value = 9.033792839530502
$wsDesc.Cells.Item(4, 3).Value = 0.4655341423056631
$wsDesc.Cells.Item(4, 4).Value = 17.61883116205967
$wsDesc.Cells.Item(4, 5).Value = 48.44270445193932
$wsDesc.Cells.Item(4, 6).Value = 22.73311184831199

# min
$wsDesc.Cells.Item(5, 2).Value = 29
$wsDesc.Cells.Item(5, 3).Value = 0
$wsDesc.Cells.Item(5, 4).Value = 94
$wsDesc.Cells.Item(5, 5).Value = 126
$wsDesc.Cells.Item(5, 6).Value = 71

# 25%
$wsDesc.Cells.Item(6, 2).Value = 48
$wsDesc.Cells.Item(6, 3).Value = 0
$wsDesc.Cells.Item(6, 4).Value = 120
$wsDesc.Cells.Item(6, 5).Value = 211
$wsDesc.Cells.Item(6, 6).Value = 134

# 50%
$wsDesc.Cells.Item(7, 2).Value = 55
$wsDesc.Cells.Item(7, 3).Value = 1
$wsDesc.Cells.Item(7, 4).Value = 130
$wsDesc.Cells.Item(7, 5).Value = 241
$wsDesc.Cells.Item(7, 6).Value = 153

# 75%
$wsDesc.Cells.Item(8, 2).Value = 61
$wsDesc.Cells.Item(8, 3).Value = 1
$wsDesc.Cells.Item(8, 4).Value = 140
$wsDesc.Cells.Item(8, 5).Value = 275
$wsDesc.Cells.Item(8, 6).Value = 166

# max
$wsDesc.Cells.Item(9, 2).Value = 77
$wsDesc.Cells.Item(9, 3).Value = 1
$wsDesc.Cells.Item(9, 4).Value = 200
$wsDesc.Cells.Item(9, 5).Value = 417
$wsDesc.Cells.Item(9, 6).Value = 202

# ---------------------------------------------------------------------
# 2. "Data Samples" sheet - header text is unchanged (cp/fbs/restecg/
#    thalach), just re-assert the values.
# ---------------------------------------------------------------------
$wsSamp.Cells.Item(1, 4).Value = "cp"
$wsSamp.Cells.Item(1, 7).Value = "fbs"
$wsSamp.Cells.Item(1, 8).Value = "restecg"
$wsSamp.Cells.Item(1, 9).Value = "thalach"

# ---------------------------------------------------------------------
# 3. "Classes Balance" sheet - updated Healthy / Sick counts.
# ---------------------------------------------------------------------
$wsBal.Cells.Item(2, 1).Value = "Healthy"
$wsBal.Cells.Item(2, 2).Value = 109
$wsBal.Cells.Item(2, 3).Value = 53
$wsBal.Cells.Item(2, 4).Value = 162

$wsBal.Cells.Item(3, 1).Value = "Sick"
$wsBal.Cells.Item(3, 2).Value = 101
$wsBal.Cells.Item(3, 3).Value = 38
$wsBal.Cells.Item(3, 4).Value = 139

# ---------------------------------------------------------------------
# 4. "Intrinsic Discrepancy" sheet - header text unchanged, values
#    updated. Values are stored as text (matching the source export),
#    so force text via a leading apostrophe.
# ---------------------------------------------------------------------
$wsDisc.Cells.Item(2, 2).Value = "age"
$wsDisc.Cells.Item(2, 3).Value = "sex"
$wsDisc.Cells.Item(2, 4).Value = "cp"
$wsDisc.Cells.Item(2, 5).Value = "restbp"
$wsDisc.Cells.Item(2, 6).Value = "chol"
$wsDisc.Cells.Item(2, 7).Value = "fbs"
$wsDisc.Cells.Item(2, 8).Value = "restecg"
$wsDisc.Cells.Item(2, 9).Value = "thalach"
$wsDisc.Cells.Item(2, 10).Value = "exang"
$wsDisc.Cells.Item(2, 11).Value = "oldpeak"
$wsDisc.Cells.Item(2, 12).Value = "slope"
$wsDisc.Cells.Item(2, 13).Value = "ca"
$wsDisc.Cells.Item(2, 14).Value = "thal"
$wsDisc.Cells.Item(2, 15).Value = "num"

$wsDisc.Cells.Item(3, 2).Value = "'0.189"
$wsDisc.Cells.Item(3, 3).Value = "'0.144"
$wsDisc.Cells.Item(3, 4).Value = "'0.593"
$wsDisc.Cells.Item(3, 5).Value = "'0.03"
$wsDisc.Cells.Item(3, 6).Value = "'0.086"
$wsDisc.Cells.Item(3, 7).Value = "'0.002"
$wsDisc.Cells.Item(3, 8).Value = "'0.067"
$wsDisc.Cells.Item(3, 9).Value = "'0.388"
$wsDisc.Cells.Item(3, 10).Value = "'0.356"
$wsDisc.Cells.Item(3, 11).Value = "'0.422"
$wsDisc.Cells.Item(3, 12).Value = "'0.319"
$wsDisc.Cells.Item(3, 13).Value = "'0.474"
$wsDisc.Cells.Item(3, 14).Value = "'0.614"
$wsDisc.Cells.Item(3, 15).Value = "'0.0"
